$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# Normalize column widths for A:F to 50 (matches existing C:F width)
$ws.Range("A:F").ColumnWidth = 49.14

# Append new test-case rows (13-18)
$ws.Range("A13").Value = '<Tier 2 User Successfully Sends International Remittance-Verify a Tier 2 (Verified) user can successfully initiate an international money transfer.>'
$ws.Range("B13").Value = 'TC-012'
$ws.Range("C13").Value = 'This test case validates the happy path for the international remittance feature as per US-401.'
$ws.Range("D13").Value = 'User must be logged in as a ''Tier 2 (Verified)'' user with sufficient funds in their account. A beneficiary must be set up.'
$ws.Range("E13").Value = '1. Log in as a Tier 2 user.
2. Navigate to the ''International Transfer'' or ''Send Money'' section.
3. Select a pre-saved beneficiary.
4. Enter the amount to send.
5. Review the exchange rate and fees.
6. Confirm the transaction.'
$ws.Range("F13").Value = 'The transfer is initiated successfully. The user sees a confirmation screen with a transaction reference number. The account balance is updated, and the transfer appears in the international transfer history with a ''Pending'' status.'

$ws.Range("A14").Value = '<Tier 1 User Blocked from International Remittance-Verify a Tier 1 (Unverified) user is not able to access or use the international remittance feature.>'
$ws.Range("B14").Value = 'TC-013'
$ws.Range("C14").Value = 'This negative test case enforces the business rule in US-401 that only verified users can send money internationally.'
$ws.Range("D14").Value = 'User must be logged in as a ''Tier 1 (Unverified)'' user.'
$ws.Range("E14").Value = '1. Log in as a Tier 1 user.
2. Attempt to navigate to the ''International Transfer'' or ''Send Money'' section.'
$ws.Range("F14").Value = 'The ''International Transfer'' option should be disabled or not visible. If accessed directly, the user should be redirected or shown a message stating they need to be a Tier 2 user to access this feature.'

$ws.Range("A15").Value = '<User Adds and Saves a New Beneficiary-Verify a user can add and save the details of a new international beneficiary.>'
$ws.Range("B15").Value = 'TC-014'
$ws.Range("C15").Value = 'This test case validates the beneficiary management feature described in US-402.'
$ws.Range("D15").Value = 'User must be logged in as a Tier 2 user.'
$ws.Range("E15").Value = '1. Log in as a Tier 2 user.
2. Navigate to the ''Beneficiaries'' or ''Recipients'' management page.
3. Click ''Add New Beneficiary''.
4. Fill in all required fields (Full Name, Country, Bank Account/Mobile Money details).
5. Save the beneficiary.'
$ws.Range("F15").Value = 'The new beneficiary is saved successfully and appears in the list of saved beneficiaries. A confirmation message is displayed.'

$ws.Range("A16").Value = '<User Deletes an Existing Beneficiary-Verify a user can delete a previously saved beneficiary.>'
$ws.Range("B16").Value = 'TC-015'
$ws.Range("C16").Value = 'This test case validates the beneficiary deletion functionality from US-402.'
$ws.Range("D16").Value = 'User must be logged in as a Tier 2 user and have at least one saved beneficiary.'
$ws.Range("E16").Value = '1. Log in as a Tier 2 user.
2. Navigate to the ''Beneficiaries'' management page.
3. Select a beneficiary from the list.
4. Click the ''Delete'' or ''Remove'' option.
5. Confirm the deletion in the confirmation prompt.'
$ws.Range("F16").Value = 'The beneficiary is successfully removed from the list of saved beneficiaries.'

$ws.Range("A17").Value = '<Verify Display of Exchange Rate and Fees Before Transfer Confirmation-Verify that the exchange rate and all applicable fees are clearly displayed to the user before they confirm an international transfer.>'
$ws.Range("B17").Value = 'TC-016'
$ws.Range("C17").Value = 'This test case validates the transparency requirement of US-403.'
$ws.Range("D17").Value = 'User is logged in as a Tier 2 user and is in the process of initiating an international transfer.'
$ws.Range("E17").Value = '1. Log in as a Tier 2 user.
2. Navigate to the ''International Transfer'' section.
3. Select a beneficiary and enter a sending amount.
4. Proceed to the confirmation/review screen.'
$ws.Range("F17").Value = 'The confirmation screen must clearly display the send amount, the exchange rate being used, the calculated recipient amount, a breakdown of all transaction fees, and the total amount to be debited. The ''Confirm'' button should be present.'

$ws.Range("A18").Value = '<Verify International Transfer History and Pagination-Verify the user can view a paginated history of their international transfers with correct statuses.>'
$ws.Range("B18").Value = 'TC-017'
$ws.Range("C18").Value = 'This test case validates the requirements of US-404 for viewing transfer history.'
$ws.Range("D18").Value = 'User must be logged in and have a history of more than 25 international transfers with various statuses (Completed, Failed, Pending).'
$ws.Range("E18").Value = '1. Log in to the platform.
2. Navigate to the ''International Transfer History'' page.
3. Observe the list of transactions and their statuses.
4. Verify the number of items on the first page.
5. Click the ''Next'' page button.'
$ws.Range("F18").Value = 'The page displays a list of past international transfers with their status. The first page shows a maximum of 25 transactions. Clicking ''Next'' loads the subsequent set of transactions.'
